# Generate Report for Handback
# Updates the handback-status report with refreshed timestamps and a
# Priority value change (ht -> mt) for the 19c20abd-... file entries.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 19c20abd-... row
$wsOverview.Range("G2").Value = "2016-09-03 02:19:05"
$wsOverview.Range("G5").Value = "2016-09-03 02:19:05"

# zh-cn sheet: Priority changes from "ht" to "mt"
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-09-03 02:18:58"
$wsZhCn.Range("H5").Value = "2016-09-03 02:18:58"
$wsZhCn.Range("K2").Value = "2016-09-03 02:19:26"
$wsZhCn.Range("K5").Value = "2016-09-03 02:19:26"

# de-de sheet: Priority changes from "ht" to "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# de-de sheet: Correspond Handoff Datetime (shares text with Overview's
# "Latest HO Xliff Generate Date") and Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-09-03 02:19:05"
$wsDeDe.Range("H5").Value = "2016-09-03 02:19:05"
$wsDeDe.Range("K2").Value = "2016-09-03 02:19:34"
$wsDeDe.Range("K5").Value = "2016-09-03 02:19:34"
